$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 57, shifting existing rows 57-58 down to 58-59.
$ws.Rows.Item(57).Insert()

# Populate the newly inserted row 57 with the new weekly entry.
$ws.Range("A57").Value = 9
$ws.Range("B57").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C57").Value = 'Metropolitana'
$ws.Range("D57").Value = 44448
$ws.Range("D57").NumberFormat = $ws.Range("D56").NumberFormat
$ws.Range("E57").Value = 13
$ws.Range("F57").Value = 100112022
$ws.Range("G57").Value = 'Arveja Verde'
$ws.Range("H57").Value = 'Perfection'
$ws.Range("I57").Value = 'Primera'
$ws.Range("J57").Value = 28
$ws.Range("K57").Value = 36000
$ws.Range("L57").Value = 37000
$ws.Range("M57").Value = 36500
$ws.Range("N57").Value = '$/malla 25 kilos'
$ws.Range("O57").Value = 'Provincia de Huasco'
$ws.Range("P57").Value = 1460
$ws.Range("Q57").Value = 25
$ws.Range("R57").Value = 'Hortaliza'
